# Workbook edit script: add a "country" column to params + two new
# power_B country rows, rename "region" -> "country" on the per-scenario
# sheets, update their "id" column values, and delete the now-redundant
# "power_B" sheet (its data lives inline in "params" from now on).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "params" sheet: insert a new "country" column (C) and populate two
#    new rows (6 & 7) describing power_B/UK and power_B/DE.
# ---------------------------------------------------------------------
$params = $wb.Worksheets("params")

$params.Columns("C").Insert()
$params.Range("C1").Value = "country"

$params.Range("A6").Value = "power_B"
$params.Range("C6").Value = "UK"
$params.Range("D6").Value = "interp"
$params.Range("E6").Value = "linear"
$params.Range("F6").Value = '{"2020-01-01":100, "2031-06-01":95}'
$params.Range("G6").Value = 1
$params.Range("H6").Value = 5
$params.Range("I6").Value = 0.05
$params.Range("J6").Value = 43617
$params.Range("K6").Value = "minute"
$params.Range("R6").Value = "x"
$params.Range("T6").Value = 16

$params.Range("A7").Value = "power_B"
$params.Range("C7").Value = "DE"
$params.Range("D7").Value = "interp"
$params.Range("E7").Value = "linear"
$params.Range("F7").Value = '{"2020-01-01":100, "2031-06-01":95}'
$params.Range("G7").Value = 2
$params.Range("H7").Value = 5
$params.Range("I7").Value = 0.05
$params.Range("J7").Value = 43617
$params.Range("K7").Value = "minute"
$params.Range("R7").Value = "x"
$params.Range("T7").Value = 17

$params.Activate()
$params.Range("A6:A7").Select()

# ---------------------------------------------------------------------
# 2. "time_A" sheet: rename header, renumber ids.
# ---------------------------------------------------------------------
$timeA = $wb.Worksheets("time_A")
$timeA.Range("A1").Value = "country"
$timeA.Range("G2").Value = 5
$timeA.Range("G3").Value = 6
$timeA.Activate()
$timeA.Range("A1").Select()

# ---------------------------------------------------------------------
# 3. "time_B" sheet: rename header, renumber ids.
# ---------------------------------------------------------------------
$timeB = $wb.Worksheets("time_B")
$timeB.Range("A1").Value = "country"
$timeB.Range("G2").Value = 7
$timeB.Range("G3").Value = 8
$timeB.Activate()
$timeB.Range("A1").Select()

# ---------------------------------------------------------------------
# 4. "power_A" sheet: rename header, renumber ids.
# ---------------------------------------------------------------------
$powerA = $wb.Worksheets("power_A")
$powerA.Range("A1").Value = "country"
$powerA.Range("G2").Value = 3
$powerA.Range("G3").Value = 15
$powerA.Activate()
$powerA.Range("A1").Select()

# ---------------------------------------------------------------------
# 5. "power_B" sheet is no longer needed as a standalone tab; its two
#    country rows now live inline in "params" (rows 6 & 7 above).
# ---------------------------------------------------------------------
$wb.Worksheets("power_B").Delete()

# ---------------------------------------------------------------------
# 6. "changes" sheet: just a selection change.
# ---------------------------------------------------------------------
$changes = $wb.Worksheets("changes")
$changes.Activate()
$changes.Range("K14").Select()

# ---------------------------------------------------------------------
# 7. "metadata" sheet: selection stays at A1.
# ---------------------------------------------------------------------
$metadata = $wb.Worksheets("metadata")
$metadata.Activate()
$metadata.Range("A1").Select()

# ---------------------------------------------------------------------
# 8. Final active tab should be "time_A".
# ---------------------------------------------------------------------
$timeA.Activate()
